$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Clear old shared-formula content first, since old A1 was a constant and
# A2:A9 held formulas referencing it. We'll overwrite wholesale below.
$ws1.Range("A1:A9").ClearContents()

# Header row
$ws1.Range("A1").Value = "Value"
$ws1.Range("B1").Value = "Formula"

# Column A: constants 1..18 in rows 2..19
for ($i = 0; $i -lt 18; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 1).Value = $i + 1
}

# Column B: formulas =A{row}*10 for rows 2..19 (rows 3..19 become one shared group)
for ($row = 2; $row -le 19; $row++) {
    $ws1.Cells.Item($row, 2).Formula = "=A" + $row + "*10"
}

$ws1.Range("B4").Select()

$wb.Worksheets.Item("Sheet2").Range("A1").Select()
$wb.Worksheets.Item("Sheet3").Range("A1").Select()
$ws1.Select()
